# Bids_template.xlsx: rename the bid-price template placeholder.
#
# The JXLS template row (row 2) drives column G with ${record.price}; the
# underlying data model field was renamed to biddingPrice, so the template
# cell must reference ${record.biddingPrice} instead. D2/E2/F2 keep their
# existing template text (salesPersonFullName / description /
# productNamesString) — only G2's text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "`${record.biddingPrice}"
